# Change the extended submission-deadline date shown on the poster from
# "June 5, 2021(15 Khordad 1400)" to "June 4, 2021(14 Khordad 1400)".
#
# The target text run "June 5, 2021(15 " is split into three runs so that
# only the day number ("5" -> "4") and the trailing "(15 " -> "(14 " text
# change while "June " keeps its own run, matching how PowerPoint would
# naturally split a run when only part of it is retyped.

$p = $ppt.ActivePresentation

$oldSnippet = "June 5, 2021(15 "
$target = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($ki = 1; $ki -le $slide.Shapes.Count; $ki++) {
        $shp = $slide.Shapes.Item($ki)

        if ($shp.Type -eq 6) {
            # msoGroup: look inside the group's member shapes too.
            for ($gi = 1; $gi -le $shp.GroupItems.Count; $gi++) {
                $inner = $shp.GroupItems.Item($gi)
                if ($inner.HasTextFrame) {
                    if ($inner.TextFrame.HasText) {
                        if ($inner.TextFrame.TextRange.Text.Contains($oldSnippet)) {
                            $target = $inner
                        }
                    }
                }
            }
        } else {
            if ($shp.HasTextFrame) {
                if ($shp.TextFrame.HasText) {
                    if ($shp.TextFrame.TextRange.Text.Contains($oldSnippet)) {
                        $target = $shp
                    }
                }
            }
        }
    }
}

if ($target -ne $null) {
    # Preserve the textbox's auto-fit size: re-typing text inside it can
    # trigger a recalculation of the box's height/width, which should stay
    # as it was since the replacement text is the same length as before.
    $origHeight = $target.Height
    $origWidth  = $target.Width

    $tr = $target.TextFrame.TextRange
    $found = $tr.Find($oldSnippet)
    $startPos = $found.Start

    # "June " (5 chars) + "5" (1 char) + ", 2021(15 " (10 chars)
    $dayRun = $tr.Characters($startPos + 5, 1)
    $dayRun.Text = "4"

    $restRun = $tr.Characters($startPos + 6, 10)
    $restRun.Text = ", 2021(14 "

    $target.Height = $origHeight
    $target.Width  = $origWidth

    Write-Host "Updated text:" $tr.Text
} else {
    Write-Host "Target text not found"
}
